# "Deleted reports and logs" — the RequiredData.xlsx test fixture was
# re-generated, which rotated the stored test credentials. Reproduce that
# content change: the "Credentials" sheet has headers in row 1
# (Username/Password) and the values in row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mngr608625"
$ws.Range("B2").Value = "uvabujU"
